$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'55.636.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.92%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.495.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +6.62%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'480.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +8.04%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'139.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +10.11%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.23%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.98%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.491.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.01%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  +6.93%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.28%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +6.03%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  +0.41%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'2.929.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +6.42%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'55.675.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.06%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "'Avalanche"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'20.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +9.18%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "'ShibaInu"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +13.45%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.492.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.52%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'4.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +10.59%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'320.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +7.39%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'10.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +9.08%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'5.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +5.42%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'57.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +4.18%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("B25").Value = "'Polygon"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.406"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +10.64%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "'Binance-PegBSC-USD"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.25%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +7.56%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'2.603.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.00%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +6.61%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0₃0780"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.18%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = "'  +0.33%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'149.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.52%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'18.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +5.17%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +8.97%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'5.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +11.39%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'3.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.63%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  +10.66%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.50%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'34.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.67%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.42%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.609"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +15.21%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +10.81%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Filecoin"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +8.13%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'Stacks"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'1.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +9.02%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  -1.27%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "'Maker"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'1.969.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.55%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "'Stellar"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0904"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +8.59%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  +5.94%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'248.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +29.39%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'4.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +11.30%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'17.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +9.76%  "
$ws.Range("E51").Style = "Normal"
